$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value  = 112
$ws.Range("I2").Value  = 336
$ws.Range("J2").Value  = 1270
$ws.Range("K2").Value  = 10
$ws.Range("L2").Value  = 355
$ws.Range("M2").Value  = 17
$ws.Range("N2").Value  = 218
$ws.Range("O2").Value  = 1
$ws.Range("P2").Value  = 5
$ws.Range("Q2").Value  = 3
$ws.Range("R2").Value  = 21
$ws.Range("T2").Value  = 222
$ws.Range("U2").Value  = 15
$ws.Range("V2").Value  = 2107
$ws.Range("X2").Value  = 2108
$ws.Range("Y2").Value  = 2
$ws.Range("Z2").Value  = 27
$ws.Range("AA2").Value = 17

$wb.Save()
